$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder question headers: B1=Ques3, C1=Ques1, D1=Ques2
$ws.Range("B1").Value = "Ques3"
$ws.Range("C1").Value = "Ques1"
$ws.Range("D1").Value = "Ques2"

# Ensure the percentage strings below are kept as literal text, not converted
# into numeric percentage values by Excel's auto-detection.
$pctRange = $ws.Range("B2:D4")
$pctRange.NumberFormat = "@"

# Row 2 - Darshan_Padia_65 - all pass (100%)
$ws.Range("B2").Value = "100.0%"
$ws.Range("C2").Value = "100.0%"
$ws.Range("D2").Value = "100.0%"

# Row 3 - Mustafa_Bharamal_78 - partial scores
$ws.Range("B3").Value = "33.3%"
$ws.Range("C3").Value = "66.7%"
$ws.Range("D3").Value = "33.3%"
$ws.Range("E3").Value = 5.333333333333333

# Row 4 - Priya_Rajani_12 - all 0%
$ws.Range("B4").Value = "0.0%"
$ws.Range("C4").Value = "0.0%"
$ws.Range("D4").Value = "0.0%"
